# "Update iterations and times" - add the three milestone headers (and
# their duration labels) above each iteration block on row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the milestone titles first, then the "(N weeks)" duration labels,
# so new shared-string entries land in the same order as the target file.
$ws.Range("F2").Value = "Mile stone 1 "
$ws.Range("J2").Value = "Mile stone 2"
$ws.Range("M2").Value = "Mile stone 3"

$ws.Range("G2").Value = "(2 weeks)"
$ws.Range("K2").Value = "(2 weeks)"
$ws.Range("N2").Value = "(3 weeks)"

# Leave the view pointed at the cell the author ended up on, un-zoomed
# back to the normal 100% view.
$ws.Range("K12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
